{"js": "// Office.js (Word JavaScript API) script\n//\n// Applies the commit \"se a\u00f1ade horarios, direccion y localidad trabajo,\n// con sus respectivas columnas de BD. check de alumnos e informes\":\n//   1. The \"Hora inicio: ${HORA_INICIO}  <tab>Hora terminaci\u00f3n: ${HORA_TERMINACION}\"\n//      line is collapsed into a single \"${HORARIOS}\" placeholder (the 16\n//      trailing spaces that followed the old closing brace are kept, since\n//      only the searched span itself is replaced).\n//   2. The \"${LOCALIDAD_EMPRESA}\" placeholder becomes \"${LOCALIDAD_TRABAJO}\".\n//   3. The \"${DIRECCION_EMPRESA}\" placeholder becomes \"${DIRECCION_TRABAJO}\".\n//\n// Note: other placeholders that merely contain the substring \"EMPRESA\"\n// (${NOMBRE_EMPRESA}, ${NOMBRE_TUTOR_EMPRESA}, ${NOMBRE_REPRESENTANTE_EMPRESA})\n// must stay untouched, so the searches below target the exact longer tokens\n// instead of a blind \"EMPRESA\" -> \"TRABAJO\" swap.\n\nconst body = context.document.body;\n\n// 1) Collapse the \"Hora inicio / Hora terminaci\u00f3n\" line into \"${HORARIOS}\"\nconst horaSearch = body.search(\"Hora inicio: ${HORA_INICIO}  \\tHora terminaci\u00f3n: ${HORA_TERMINACION}\", {\n  matchCase: true\n});\nhoraSearch.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < horaSearch.items.length; i++) {\n  horaSearch.items[i].insertText(\"${HORARIOS}\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) LOCALIDAD_EMPRESA -> LOCALIDAD_TRABAJO\nconst localidadSearch = body.search(\"LOCALIDAD_EMPRESA\", { matchCase: true });\nlocalidadSearch.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < localidadSearch.items.length; i++) {\n  localidadSearch.items[i].insertText(\"LOCALIDAD_TRABAJO\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 3) DIRECCION_EMPRESA -> DIRECCION_TRABAJO\nconst direccionSearch = body.search(\"DIRECCION_EMPRESA\", { matchCase: true });\ndireccionSearch.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < direccionSearch.items.length; i++) {\n  direccionSearch.items[i].insertText(\"DIRECCION_TRABAJO\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Word COM interop script\n#\n# Applies the commit \"se a\u00f1ade horarios, direccion y localidad trabajo,\n# con sus respectivas columnas de BD. check de alumnos e informes\":\n#   1. The \"Hora inicio: ${HORA_INICIO}  <tab>Hora terminaci\u00f3n: ${HORA_TERMINACION}\"\n#      line is collapsed into a single \"${HORARIOS}\" placeholder (the 16\n#      trailing spaces after the old closing brace are kept, since Find\n#      only touches the matched span).\n#   2. The \"${LOCALIDAD_EMPRESA}\" placeholder becomes \"${LOCALIDAD_TRABAJO}\".\n#   3. The \"${DIRECCION_EMPRESA}\" placeholder becomes \"${DIRECCION_TRABAJO}\".\n#\n# Note: other placeholders that merely contain the substring \"EMPRESA\"\n# (${NOMBRE_EMPRESA}, ${NOMBRE_TUTOR_EMPRESA}, ${NOMBRE_REPRESENTANTE_EMPRESA})\n# must stay untouched, so the searches below target the exact longer tokens\n# instead of a blind \"EMPRESA\" -> \"TRABAJO\" swap.\n# `$` is escaped as `` `$ `` and the tab character as `` `t `` so the\n# double-quoted strings are not treated as PowerShell variable expansions.\n\n$d = $word.ActiveDocument\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n# 1) Collapse the \"Hora inicio / Hora terminaci\u00f3n\" line into \"${HORARIOS}\"\n$horaSearch = \"Hora inicio: `${HORA_INICIO}  `tHora terminaci\u00f3n: `${HORA_TERMINACION}\"\n$range1 = $d.Content\n$range1.Find.ClearFormatting()\n$range1.Find.Execute($horaSearch, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"`${HORARIOS}\", $wdReplaceAll)\n\n# 2) LOCALIDAD_EMPRESA -> LOCALIDAD_TRABAJO\n$range2 = $d.Content\n$range2.Find.ClearFormatting()\n$range2.Find.Execute(\"LOCALIDAD_EMPRESA\", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"LOCALIDAD_TRABAJO\", $wdReplaceAll)\n\n# 3) DIRECCION_EMPRESA -> DIRECCION_TRABAJO\n$range3 = $d.Content\n$range3.Find.ClearFormatting()\n$range3.Find.Execute(\"DIRECCION_EMPRESA\", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"DIRECCION_TRABAJO\", $wdReplaceAll)\n"}
